# Generate Report for Handback
# Updates the handback-status report timestamps and status ("ht" -> "mt")
# on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 16:17:13"
$wsOverview.Range("G3").Value = "2016-08-22 16:17:13"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-22 16:17:00"
$wsZhCn.Range("H3").Value = "2016-08-22 16:17:00"
$wsZhCn.Range("K2").Value = "2016-08-22 16:17:30"
$wsZhCn.Range("K3").Value = "2016-08-22 16:17:30"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-22 16:17:13"
$wsDeDe.Range("H3").Value = "2016-08-22 16:17:13"
$wsDeDe.Range("K2").Value = "2016-08-22 16:17:37"
$wsDeDe.Range("K3").Value = "2016-08-22 16:17:37"
